$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows of data (row 12 and row 13), following the same
# json_path_before_change / json_path_after_change / instructions layout
# used by the existing rows.
$ws.Range("A12").Value = "jsons_train/rapport_original.json"
$ws.Range("B12").Value = 'Desole je ne reconnais pas le terme "xljfsn"'
$ws.Range("C12").Value = "Ajoute un xljfsn du nombre distinct de session_id"

$ws.Range("A13").Value = "jsons_train/rapport_original.json"
$ws.Range("B13").Value = 'Desole mais le visuel "graphique en etoile" n''est pas disponible dans Power BI'
$ws.Range("C13").Value = "Ajoute un graphique en etoile des variables de la table"

# Resize / re-fit the columns now that a third column is in play.
$ws.Columns.Item(1).ColumnWidth = 31.893229166666668
$ws.Columns.Item(2).ColumnWidth = 60.072916666666664
$ws.Columns.Item(3).ColumnWidth = 81.34635416666667

# Move the active selection, matching the post-edit saved state.
$ws.Range("B7").Select()
